# Fruta / hortaliza, semanal
# Insert a new weekly data row (2021-09-13, bandeja 4 kilos, Brasil) before
# the current row 32, shifting the existing rows 32-47 down to 33-48.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new row at position 32; existing rows 32:47 shift down to 33:48.
$ws.Rows.Item(32).Insert()

# Copy the formatting of the row immediately below (old row 32, now row 33)
# onto the freshly inserted row so styles (e.g. the date format on column D)
# carry over, without touching the whole 1..16384 column span.
$ws.Range("A33:T33").Copy()
$ws.Range("A32:T32").PasteSpecial(-4122)  # xlPasteFormats
$excel.CutCopyMode = $false

$ws.Cells.Item(32, 1).Value = 7
$ws.Cells.Item(32, 2).Value = "Terminal Hortofrutícola Agro Chillán"
$ws.Cells.Item(32, 3).Value = "Ñuble"
$ws.Cells.Item(32, 4).Value = "2021-09-13"
$ws.Cells.Item(32, 5).Value = 16
$ws.Cells.Item(32, 6).Value = "Fruta"
$ws.Cells.Item(32, 7).Value = 100108
$ws.Cells.Item(32, 8).Value = "Tropicales y subtropicales"
$ws.Cells.Item(32, 9).Value = 100108002
$ws.Cells.Item(32, 10).Value = "Mango"
$ws.Cells.Item(32, 11).Value = "Sin especificar"
$ws.Cells.Item(32, 12).Value = "Primera"
$ws.Cells.Item(32, 13).Value = 60
$ws.Cells.Item(32, 14).Value = 8500
$ws.Cells.Item(32, 15).Value = 9000
$ws.Cells.Item(32, 16).Value = 8750
$ws.Cells.Item(32, 17).Value = '$/bandeja 4 kilos'
$ws.Cells.Item(32, 18).Value = "Brasil"
$ws.Cells.Item(32, 19).Value = 2188
$ws.Cells.Item(32, 20).Value = 4
